# feature : added feature to get the active customer in the customer landing page
#
# The "User Data" sheet is a simple users table (User ID, Email, First Name,
# Last Name, Address, User Type). Getting the active customer on the
# customer landing page means a new customer row shows up in this table -
# add the next user (id 3) as an active "Customer" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Data")

$row = 4
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "sursau@outlook.com"
$ws.Cells.Item($row, 3).Value = "Saurav"
$ws.Cells.Item($row, 4).Value = "S"
$ws.Cells.Item($row, 5).Value = "uk"
$ws.Cells.Item($row, 6).Value = "Customer"
